$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-14 Friday" "2025-11-15 Saturday"

Replace-Text "753÷4=" "786÷5="
Replace-Text "804÷5=" "920÷7="
Replace-Text "541÷5=" "468÷3="
Replace-Text "202÷4=" "596÷4="
Replace-Text "811÷6=" "519÷8="
Replace-Text "870÷4=" "647÷6="
Replace-Text "802÷4=" "824÷6="
Replace-Text "434÷2=" "449÷9="
Replace-Text "967÷9=" "271÷8="
Replace-Text "197÷6=" "170÷9="
Replace-Text "575÷8=" "580÷3="
Replace-Text "469÷4=" "898÷7="
Replace-Text "564÷2=" "172÷8="
Replace-Text "833÷9=" "276÷9="
Replace-Text "848÷6=" "856÷4="
Replace-Text "216÷3=" "357÷6="
Replace-Text "234÷8=" "134÷6="
Replace-Text "907÷7=" "808÷2="
Replace-Text "254÷2=" "474÷6="
Replace-Text "438÷4=" "858÷7="
Replace-Text "749÷3=" "400÷9="
Replace-Text "668÷4=" "765÷8="
Replace-Text "586÷5=" "527÷4="
Replace-Text "235÷5=" "905÷4="
Replace-Text "629÷9=" "488÷7="
